$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "transfer skip and update add new character event":
# Insert a new row above the last row (old row 25, the "skip"/Empty event),
# pushing it down to row 26, and use the freed row 25 for the
# "add new character" event.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the "add 1005 new character" event data
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "add 1005 new character"
$ws.Range("C25").Value = "add 1005 new character"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = "null"
$ws.Range("F25").Value = "[40]"

# Reset the view back to the top-left of the sheet (clears the stored
# topLeftCell scroll position left over from A10)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Keep the selection on E25, matching the saved selection in the sheet
$ws.Range("E25").Select() | Out-Null
